{"js": "// 1) Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\nconst dateMatches = context.document.body.search(\"September 19, 2025\", { matchCase: true });\ndateMatches.load(\"items\");\nawait context.sync();\nif (dateMatches.items.length > 0) {\n  dateMatches.items[0].insertText(\"September 21, 2025\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Split the mailing-address paragraph \"999 Story Road, San Jose CA 95122\"\n//    (the copy that lives directly in the body, NOT the duplicate inside the\n//    PROPERTY ADDRESS table) into two paragraphs:\n//      \"999 Story Road\"\n//      \"San Jose, CA 95122\"\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Gather every paragraph whose text matches, together with a flag telling us\n// whether it lives inside a table cell (the PROPERTY ADDRESS table also has\n// this exact text and must stay untouched).\nconst candidates = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text === \"999 Story Road, San Jose CA 95122\") {\n    const cell = p.parentTableCellOrNullObject;\n    cell.load(\"isNullObject\");\n    candidates.push({ paragraph: p, cell: cell });\n  }\n}\nawait context.sync();\n\nlet addressParagraph = null;\nfor (const candidate of candidates) {\n  if (candidate.cell.isNullObject) {\n    addressParagraph = candidate.paragraph;\n    break;\n  }\n}\n\nif (addressParagraph) {\n  addressParagraph.insertText(\"999 Story Road\", \"Replace\");\n  addressParagraph.insertParagraph(\"San Jose, CA 95122\", \"After\");\n  await context.sync();\n}\n\n// 3) Remove the extra empty \"No Spacing\" paragraph that sits right after the\n//    \"Board of Directors\" paragraph in the signature block.\nconst allParagraphs = context.document.body.paragraphs;\nallParagraphs.load(\"items/text\");\nawait context.sync();\n\nlet boardParagraph = null;\nfor (let i = 0; i < allParagraphs.items.length; i++) {\n  if (allParagraphs.items[i].text.indexOf(\"Board of Directors\") !== -1) {\n    boardParagraph = allParagraphs.items[i];\n    break;\n  }\n}\n\nif (boardParagraph) {\n  const nextParagraph = boardParagraph.getNext();\n  nextParagraph.load(\"text,style\");\n  await context.sync();\n  if (nextParagraph.text === \"\" && nextParagraph.style === \"No Spacing\") {\n    nextParagraph.delete();\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\n$find = $d.Content\n$find.Find.ClearFormatting()\n$find.Find.Text = \"September 19, 2025\"\n$find.Find.Forward = $true\n$dateFound = $find.Find.Execute()\nif ($dateFound) {\n    $find.Text = \"September 21, 2025\"\n}\n\n# 2) Split the mailing-address paragraph \"999 Story Road, San Jose CA 95122\"\n#    (the copy that lives directly in the body, NOT the one inside the\n#    PROPERTY ADDRESS table) into two paragraphs:\n#      \"999 Story Road\"\n#      \"San Jose, CA 95122\"\n$wdWithInTable = 12\n$addressParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    $paraText = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($paraText -eq \"999 Story Road, San Jose CA 95122\") {\n        if ($p.Range.Information($wdWithInTable) -eq $false) {\n            $addressParagraph = $p\n            break\n        }\n    }\n}\n\nif ($addressParagraph -ne $null) {\n    $r = $addressParagraph.Range\n    $r.Text = \"999 Story Road\"\n    $r.InsertParagraphAfter()\n    $newParagraph = $addressParagraph.Next()\n    $newParagraph.Range.Text = \"San Jose, CA 95122\"\n}\n\n# 3) Remove the extra empty \"No Spacing\" paragraph that sits right after the\n#    \"Board of Directors\" paragraph in the signature block.\n$boardParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Board of Directors*\") {\n        $boardParagraph = $p\n    }\n}\n\nif ($boardParagraph -ne $null) {\n    $nextParagraph = $boardParagraph.Next()\n    if ($nextParagraph -ne $null -and $nextParagraph.Style.NameLocal -eq \"No Spacing\" -and $nextParagraph.Range.Text.Trim() -eq \"\") {\n        $nextParagraph.Range.Delete()\n    }\n}\n"}
